# Apply cryptocurrency price/volume updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '58.074.36'
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +1.89%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '2.359.43'
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +1.76%  '
$ws.Cells.Item(4, 5).Value = '  -0.34%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '541.60'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +2.30%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '136.18'
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.53%  '
$ws.Cells.Item(8, 5).Value = '  +5.29%  '
$ws.Cells.Item(9, 5).Value = '  +1.31%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '5.58'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +4.01%  '
$ws.Cells.Item(11, 5).Value = '  -0.74%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.356'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.51%  '
$ws.Cells.Item(13, 5).Value = '  +2.86%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '2.782.83'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.50%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '58.088.46'
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.81%  '
$ws.Cells.Item(16, 5).Value = '  +1.61%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '2.362.85'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.70%  '
$ws.Cells.Item(18, 5).Value = '  +3.35%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '333.32'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.98%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '4.28'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +2.36%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '6.80'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.32%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.13%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '62.78'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.93%  '
$ws.Cells.Item(24, 5).Value = '  +0.20%  '
$ws.Cells.Item(25, 5).Value = '  -2.33%  '
$ws.Cells.Item(26, 5).Value = '  +0.72%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '1.38'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +2.43%  '
$ws.Cells.Item(28, 2).Value = 'Monero'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '172.92'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -0.21%  '
$ws.Cells.Item(29, 2).Value = 'PancakeSwap'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '1.75'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +2.01%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0740'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +2.15%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '6.17'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.97%  '
$ws.Cells.Item(32, 5).Value = '  +11.61%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '18.55'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.36%  '
$ws.Cells.Item(34, 5).Value = '  +0.04%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '4.24'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +6.46%  '
$ws.Cells.Item(36, 5).Value = '  +0.77%  '
$ws.Cells.Item(37, 5).Value = '  +0.40%  '
$ws.Cells.Item(38, 5).Value = '  +3.77%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '39.39'
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '145.83'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -2.39%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '294.28'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +3.73%  '
$ws.Cells.Item(42, 5).Value = '  +1.39%  '
$ws.Cells.Item(43, 5).Value = '  +1.50%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.0950'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +2.03%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '19.24'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +2.26%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.0503'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.65%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '0.564'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +1.07%  '
$ws.Cells.Item(48, 5).Value = '  +2.82%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '17.53'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +0.28%  '
$ws.Cells.Item(50, 2).Value = 'Polygon'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.382'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.00%  '
$ws.Cells.Item(51, 5).Value = '  +0.37%  '
